# Append one new data row (row 69) to Sheet1, mirroring the existing rows:
#   A: date string "2025/10/06" (stored as text, NOT auto-converted to a date)
#   B: weekday string "月"
#   C: hour number 16
#   D: ranking number 6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 69

# Format column A as Text first so Excel doesn't auto-convert the
# "yyyy/mm/dd"-looking string into a real date serial number, then drop the
# number format back to General (matching the rest of the sheet, which has
# no explicit per-cell number format) once the text value is in place.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025/10/06"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "月"
$ws.Cells.Item($row, 3).Value = 16
$ws.Cells.Item($row, 4).Value = 6
